# Appends 10 new arrival rows (rows 126-135) to the "Main Data" sheet,
# continuing the flight-arrivals log for "Friday, Jan 13".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=NUMBER, B=DATE, C=TIME, D=FLIGHT, E=FROM, F=SHORT,
#          G=AIRLINE, H=MODEL, I=AIRCFAT ID, J=STATUS, K=(blank/DIFFERENCE header),
#          L=difference text, M=(blank)
$rows = @(
    @(125, "Friday, Jan 13", "7:00 AM",  "FR3285", "Riga",      "(RIX)", "Ryanair ", "B738", "(SP-RKV)", "6:41 AM",  "0 hours, -19 minutes"),
    @(126, "Friday, Jan 13", "9:05 AM",  "FR2019", "London",    "(STN)", "Ryanair ", "B738", "(EI-EVP)", "8:47 AM",  "0 hours, -18 minutes"),
    @(127, "Friday, Jan 13", "9:20 AM",  "FR1933", "Leeds",     "(LBA)", "Ryanair ", "B738", "(EI-EKY)", "9:48 AM",  "0 hours, 28 minutes"),
    @(128, "Friday, Jan 13", "9:40 AM",  "FR6944", "Barcelona", "(BCN)", "Ryanair ", "B738", "(EI-DYC)", "9:20 AM",  "0 hours, -20 minutes"),
    @(129, "Friday, Jan 13", "10:20 AM", "FR4529", "Oslo",      "(TRF)", "Ryanair ", "B738", "(SP-RKW)", "9:57 AM",  "0 hours, -23 minutes"),
    @(130, "Friday, Jan 13", "11:30 AM", "FR1573", "Vienna",    "(VIE)", "Ryanair ", "B738", "(SP-RKV)", "11:24 AM", "0 hours, -6 minutes"),
    @(131, "Friday, Jan 13", "11:50 AM", "FR1055", "Brussels",  "(CRL)", "Ryanair ", "B738", "(SP-RSP)", "11:40 AM", "0 hours, -10 minutes"),
    @(132, "Friday, Jan 13", "12:05 PM", "FR1973", "Catania",   "(CTA)", "Ryanair ", "B738", "(SP-RKD)", "12:07 PM", "0 hours, 2 minutes"),
    @(133, "Friday, Jan 13", "12:10 PM", "FR1889", "Paris",     "(BVA)", "Buzz ",    "B38M", "(SP-RZG)", "12:01 PM", "0 hours, -9 minutes"),
    @(134, "Friday, Jan 13", "12:35 PM", "FR2027", "Dublin",    "(DUB)", "Buzz ",    "B38M", "(SP-RZE)", "12:23 PM", "0 hours, -12 minutes")
)

$startRow = 126
foreach ($row in $rows) {
    $r = $startRow
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 12).Value = $row[10]
    $startRow = $startRow + 1
}
